$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("userData")

$ws.Range("A6").Value = "lakshmi"
$ws.Range("B6").Value = "lakshmi12@gmail.com"
$ws.Range("C6").Value = "lakshmi12345"
$ws.Range("D6").Value = 44

$ws.Hyperlinks.Add($ws.Range("B6"), "mailto:lakshmi12@gmail.com")
$ws.Range("B6").Style = "Hyperlink"

$ws.Range("C6").Select()
